# Updating filtered feeds from workflow
#
# Two brand-new rows are appended for the Myriad Genetics / Sophia Genetics
# companion diagnostic story (one row per source link: genomeweb + 360dx).
# (Existing row 49's shared-string indices shift in the underlying XML only
# because two new strings are inserted earlier in the table - its actual
# link/keyword/title content is unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 50: GenomeWeb link ------------------------------------------
$linkText50 = "https://www.genomeweb.com/cancer/myriad-genetics-sophia-genetics-collaborate-cancer-liquid-biopsy-companion-diagnostic"
$ws.Range("A50").Value2 = $linkText50
$ws.Range("B50").Value2 = "companion diagnostic"
$ws.Range("C50").Value2 = "Myriad Genetics, Sophia Genetics to Collaborate on Cancer Liquid Biopsy Companion Diagnostic"
$ws.Hyperlinks.Add($ws.Range("A50"), $linkText50) | Out-Null
$ws.Range("A50").Style = "Hyperlink"

# --- New row 51: 360Dx link ------------------------------------------
$linkText51 = "https://www.360dx.com/cancer/myriad-genetics-sophia-genetics-collaborate-cancer-liquid-biopsy-companion-diagnostic"
$ws.Range("A51").Value2 = $linkText51
$ws.Range("B51").Value2 = "companion diagnostic"
$ws.Range("C51").Value2 = "Myriad Genetics, Sophia Genetics to Collaborate on Cancer Liquid Biopsy Companion Diagnostic"
$ws.Hyperlinks.Add($ws.Range("A51"), $linkText51) | Out-Null
$ws.Range("A51").Style = "Hyperlink"
